$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: 20240304
# Write the date as a text formula first (="20240304" evaluates to the text
# string "20240304"), convert that formula result to a plain value in place,
# then stamp the same number format / font / border / alignment as the other
# date cells (A2) onto it by copying formats only - this reuses the existing
# style record instead of minting a new one.
$ws.Range("A7").Formula = '="20240304"'
$ws.Range("A7").Copy()
$ws.Range("A7").PasteSpecial(-4163)
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B7").Value = 72
$ws.Range("C7").Value = 156
$ws.Range("D7").Value = 69
$ws.Range("E7").Value = 107
$ws.Range("F7").Value = 162
$ws.Range("G7").Value = 63

# Row 8: 20240305
$ws.Range("A8").Formula = '="20240305"'
$ws.Range("A8").Copy()
$ws.Range("A8").PasteSpecial(-4163)
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("B8").Value = 71
$ws.Range("C8").Value = 339
$ws.Range("D8").Value = 68
$ws.Range("E8").Value = 107
$ws.Range("F8").Value = 167
$ws.Range("G8").Value = 58

$excel.CutCopyMode = $false
